$d = $word.ActiveDocument

$replacements = @(
    @("0.997 and coefficients", "0.9642 and coefficients"),
    @("3.0137, and", "3.0004, and"),
    @("1.9574", "1.9829"),
    @("0.9974 and coefficients", "0.9635 and coefficients"),
    @("3.0123, and", "2.9995, and"),
    @("1.9562", "1.9828"),
    @("array([0.97999721]) and coefficients", "array([0.96075922]) and coefficients"),
    @("array([3.00639323]), and", "array([2.9797648]), and"),
    @("array([1.96080487])", "array([1.98935137])")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
